$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "magapoke_2026-02-18"

$newSheet.Cells.Item(1, 1).Value = "rank"
$newSheet.Cells.Item(1, 2).Value = "title"

$titles = @(
    '怨霊日和',
    'ハンドレッドノート－アグリーダック－',
    'ハードワーカー中田',
    'ペンの夢に紅をさす',
    '黒月のイェルクナハト',
    'ドリーム☆ジャンボ☆ガール',
    '黄昏町プリズナーズ',
    'アイドラトリィ',
    'K-9~警視庁公安部公安第9課異能対策係~',
    '【爆アド】生まれた直後から最強悪霊と脳内バトルしてたら魔力量が測定可能域を超えてました〜悪憑の子の謙虚な覇道〜',
    'ゼロとヒャク',
    'せいぶつ部の田辺くん',
    'ともだちづくり',
    '追放されなかった男　～二度目の人生は土下座から始まりました～',
    '普通の本はありません！',
    '屋根の下のアルテミス',
    'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜',
    'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！',
    '眠れる森のレガ',
    'お前がヤったんだろ！',
    'あの島の海音荘',
    '春くらり',
    '白鳥運子は31画',
    'GURU',
    'MYS',
    'わが投資術　市場は誰に微笑むか',
    '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～',
    '君が監督！',
    '人生逆転ダンジョン',
    'その青春',
    '邪目さんは邪神です',
    '白銀のキュイジーヌ～明治外交官の料理人～',
    '鳴るさんだぁ',
    'ハプスブルク家の華麗なる受難',
    '平成転生',
    '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～',
    '篝家の８兄弟',
    '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～',
    '歪みの虜',
    'ch登録お願いします！',
    'JK Biker',
    'イエティ、とある日々',
    'ナキナギ',
    '明智ナンバーワン',
    '宇曽田みのりの代用料理',
    '花子狩り',
    'ナマイキ旭ちゃんをわからせたい',
    'じゅーくぼっくす',
    '永久のユウグレ',
    'きゃわるり方程式',
    '夜鐘のキト'
)

for ($i = 0; $i -lt $titles.Length; $i++) {
    $r = $i + 2
    $newSheet.Cells.Item($r, 1).Value = $i + 1
    $newSheet.Cells.Item($r, 2).Value = $titles[$i]
}
